$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Logs" sheet: append row 6 with the new test-mail entry
# ---------------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A6").Value = "Wil je 100 stuks M5-bouten bestellen?"
$logs.Range("B6").Value = "mailmind.test@zohomail.eu"
$logs.Range("C6").Value = "Testmail #4: Wil je 100 stuks M5-bouten bestellen?"
$logs.Range("D6").Value = "Bestelling / Levering"
$logs.Range("E6").Value = "Beste afzender,`nBedankt voor je interesse in het bestellen van M5-bouten. Helaas kan ik als e-mailassistent geen bestellingen plaatsen. Voor het bestellen van producten kun je terecht op onze website of contact opnemen met onze verkoopafdeling.`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$logs.Range("F6").Value = "2025-07-27 19:18:59"
$logs.Range("G6").Value = "Ja"
$logs.Range("H6").Value = "Nee"
$logs.Range("I6").Value = "Ja"
$logs.Range("J6").Value = "Nee"

# Writing the multi-line "Antwoord" text into a brand-new row makes the
# engine stamp an explicit autofit row height (ht=/customHeight=1); re-running
# AutoFit drops that explicit height again so row 6 matches the other,
# height-less rows.
$logs.Rows.Item(6).AutoFit()

# Extend the conditional-formatting ranges (previously row 2-5) so they also
# cover the freshly added row 6, one ModifyAppliesToRange per rule group.
$logs.Range("D2:D5").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D6"))
$logs.Range("G2:G5").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G6"))
$logs.Range("H2:H5").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H6"))
$logs.Range("I2:I5").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I6"))
$logs.Range("J2:J5").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J6"))

# ---------------------------------------------------------------------------
# 2. "Dashboard" sheet: append row 4 with the new category tally
# ---------------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A4").Value = "Bestelling / Levering"
$dash.Range("B4").Value = 1

# ---------------------------------------------------------------------------
# 3. Chart on the Dashboard sheet: extend category/value series to row 4
# ---------------------------------------------------------------------------
$chart = $dash.ChartObjects().Item(1).Chart
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$4,'Dashboard'!`$B`$2:`$B`$4,1)"
